$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest Lotofacil draw results to append: Concurso (A), Bola1..Bola15 (B..P)
$data = @(
    @(3541,1,2,3,4,6,8,13,14,16,17,18,19,21,23,25),
    @(3542,1,2,5,8,11,12,13,15,18,19,21,22,23,24,25),
    @(3543,1,4,5,6,7,8,9,14,15,16,17,20,21,22,25),
    @(3544,1,4,5,6,8,10,11,12,15,17,20,21,22,23,25),
    @(3545,1,4,5,6,7,8,11,12,14,17,19,20,21,22,25),
    @(3546,2,3,5,7,8,11,12,13,18,19,20,22,23,24,25),
    @(3547,2,5,6,7,9,10,14,15,17,18,19,20,21,24,25),
    @(3548,4,8,9,10,11,12,13,15,16,17,19,21,23,24,25),
    @(3549,2,5,7,8,9,11,12,14,17,20,21,22,23,24,25),
    @(3550,1,4,6,7,8,12,13,15,16,18,19,20,22,23,24),
    @(3551,1,2,3,7,8,10,11,12,15,16,17,19,21,22,25)
)

$startRow = 471
$numRows = $data.Count
$numCols = $data[0].Count

# Build a 2D array so the whole block can be written in one Range.Value call.
$arr = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$lastRow = $startRow + $numRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($lastRow, 16))
$rng.Value = $arr

# Select the freshly-entered block, matching where the workbook was left
# after typing in the new results (B471:P481, with B471 active).
$selRange = $ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($lastRow, 16))
$selRange.Select()
